# Insert a new weekly record at row 312 of the data set, shifting the
# existing rows 312-370 down to 313-371 (dimension grows from R370 to R371).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 312; Excel shifts rows
# 312..370 down to 313..371 and carries their formatting with them.
$ws.Rows(312).Insert()

# Populate the newly inserted row 312 with the new weekly observation.
$ws.Range("A312").Value2 = 9
$ws.Range("B312").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C312").Value2 = "Metropolitana"
$ws.Range("D312").Value2 = 44637
$ws.Range("E312").Value2 = 13
$ws.Range("F312").Value2 = 100112039
$ws.Range("G312").Value2 = "Ciboulette"
$ws.Range("H312").Value2 = "Sin especificar"
$ws.Range("I312").Value2 = "Primera"
$ws.Range("J312").Value2 = 550
$ws.Range("K312").Value2 = 1000
$ws.Range("L312").Value2 = 1500
$ws.Range("M312").Value2 = 1273
$ws.Range("N312").Value2 = "`$/docena de atados"
$ws.Range("O312").Value2 = "Provincia de Chacabuco"
$ws.Range("P312").Value2 = 424
$ws.Range("Q312").Value2 = 3
$ws.Range("R312").Value2 = "Hortaliza"
